# Adds a new row (row 9) to the "数组" (Array) sheet describing LeetCode
# #169 "Majority Element", solved with the Boyer-Moore majority vote
# algorithm ("摩尔投票法"), matching the commit
# "majority moer denote with array".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")
$ws.Activate()

$problem = "给定一个大小为 n 的数组，找到其中的多数元素。多数元素是指在数组中出现次数大于 ⌊ n/2 ⌋ 的元素。 `n 你可以假设数组是非空的，并且给定的数组总是存在多数元素。 `n 示例 1:`n 输入: [3,2,3]`n输出: 3 `n 示例 2: `n 输入: [2,2,1,1,1,2,2]`n输出: 2`n Related Topics 位运算 数组 分治算法"

$solution = "1 题目要求一定存在这样一个元素；只会有一个元素是多数元素`n2 摩尔投票法https://leetcode-cn.com/problems/majority-element/solution/tu-jie-mo-er-tou-piao-fa-python-go-by-jalan/`n3 设定第一个元素num出现次数count是1`n4 如果上一个元素不同于当前元素,count--`n5 如果count=0，当前元素作为基准元素，count=1`n6 迭代结束，基准元素就是多数元素"

$keyword = "摩尔投票法"

# No. / leetcode number
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 169
# 解题方法 / 题目 / 解题关键词 (filled in this order so new shared-string
# entries land at the same indices the original author produced: solution
# text first, then the problem statement, then the keyword)
$ws.Cells.Item(9, 4).Value = $solution
$ws.Cells.Item(9, 3).Value = $problem
$ws.Cells.Item(9, 5).Value = $keyword
# 时间复杂度 / 空间复杂度 (same as the other array rows: O(N) / O(1))
$ws.Cells.Item(9, 6).Value = $ws.Cells.Item(8, 6).Value()
$ws.Cells.Item(9, 7).Value = $ws.Cells.Item(8, 7).Value()

# Match the row height used for similarly sized entries in this sheet.
$ws.Rows.Item(9).RowHeight = 352

# Scroll/select so the newly added row is in view, same as the author did.
$ws.Range("D9").Select()
